{"js": "// 1. Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) left over from\n//    the previous edit session.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Collapse the \">>>  your stuff after this line >>>\" paragraph (which is\n//    split across several runs, with proofing-error markers in between) into\n//    a single run reading \">>>  your stuff after this line >>>\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"your stuff after this line\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (targetParagraph) {\n  targetParagraph.clear();\n  await context.sync();\n  targetParagraph.insertText(\">>>  your stuff after this line >>>\", \"Start\");\n  await context.sync();\n}\n\n// 3. Fill the final empty paragraph with Sarah's comment for Assignment 1.\nparagraphs.load(\"text\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\n  \"This is Sarah making a change for the purpose of Assignment 1.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the leftover \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) from\n#    the heading line.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Collapse the \">>>  your stuff after this line >>>\" paragraph (currently\n#    split across several runs with proofing-error markers) into a single run.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \">>>  your stuff after this line >>>\"\n$find.Replacement.Text = \">>>  your stuff after this line >>>\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 3. Fill the final empty paragraph with Sarah's comment for Assignment 1.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.Text = \"This is Sarah making a change for the purpose of Assignment 1.\"\n"}
